# Added course and permissive static filters
# Rename header labels on row 1 of the "template" sheet:
#   "Geofence Name"  -> "Zone"      (column F)
#   "Activity Type"  -> "Activity"  (column E)
#   "Maintenance Name" -> "Maintenance" (column G)
# and move the active selection from H10 to H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Zone"
$ws.Range("E1").Value = "Activity"
$ws.Range("G1").Value = "Maintenance"

$ws.Range("H2").Select()
